# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values per row/column (B,C,D,E,G) - F is unchanged.
$data = @{
    2 = @{ B = 0.3048080303191223; C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732;  G = 2.626907116734944 }
    3 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144 }
    4 = @{ B = 0.01514828764759746; C = 0.04240448674262143; D = 3.900430680208489; E = 0.496779210170732;  G = 4.45476266476944 }
    5 = @{ B = 0.3048080303191223; C = 1.667794583268128;  D = 3.900430680208489; E = 8.660232485948974;  G = 14.53326577974471 }
    6 = @{ B = 0.3048080303191223; C = 0.04240448674262143; D = 0.8054896365839992; E = 8.660232485948974; G = 9.812934639594717 }
    7 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732;  G = 6.201049113329182 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
